$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Sr.No"
$ws.Range("B1").Value = "Name"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Devu"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Hina"

$ws.Range("B3").Select()
